# Edit the "nomenclature" worksheet's regex column (column B).
# Most of the regex strings there had a trailing "$" removed from the
# end of the pattern; one entry (CD4, row 31) instead had "\s" inserted
# before the trailing ".*$".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nomenclature")

$ws.Range("B8").Value  = ".*Alexa\s*-*Fluor\s*-*488\s*-*A.*"
$ws.Range("B9").Value  = ".*Alexa\s*-*Fluor\s*-*647\s*-*A.*"
$ws.Range("B10").Value = ".*APC\s*-*Cy7\s*-*A.*"
$ws.Range("B11").Value = ".*Alexa\s*-*Fluor\s*-*405\s*-*A.*"
$ws.Range("B12").Value = ".*AmCyan\s*-*A.*"
$ws.Range("B13").Value = ".*BV605\s*-*A.*"
$ws.Range("B14").Value = ".*BV711\s*-*A.*"
$ws.Range("B15").Value = ".*PE\s*-*A.*"
$ws.Range("B16").Value = ".*PE\s*-*CF594\s*-*A.*"
$ws.Range("B17").Value = ".*7\s*-*AAD\s*-*A.*"
$ws.Range("B18").Value = ".*PE\s*-*Cy5\s*-*\.*5\s*-*A.*"
$ws.Range("B19").Value = ".*PE\s*-*Cy\s*-*7\s*-*A.*"

$ws.Range("B21").Value = ".*CXCR\s*-*3.*"
$ws.Range("B22").Value = ".*CD\s*-*161.*"
$ws.Range("B23").Value = ".*CD\s*-*3.*"
$ws.Range("B24").Value = ".*CCR\s*-*7.*"
$ws.Range("B25").Value = ".*L\/*D.*"
$ws.Range("B26").Value = ".*VA\s*-*7\.*-*2.*"
$ws.Range("B27").Value = ".*CD\s*-*8.*"
$ws.Range("B28").Value = ".*v\s*d[elta]*\s*-*2.*"
$ws.Range("B29").Value = ".*CD\s*-*45\s*-*RA.*"
$ws.Range("B30").Value = ".*PAN\s*-*GD.*"
$ws.Range("B31").Value = ".*CD\s*-*4\s.*$"
$ws.Range("B32").Value = ".*CD\s*-*27.*"

$ws.Range("B29").Select()
